$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Footer (first-page), holds the Pearson logo currently named "image2.png",
#     docPr/cNvPr id="3" -> rename to "image1.png" ---
$ftrFirst = $sec.Footers(2)
$xml = $ftrFirst.Range.WordOpenXML
$xml = $xml.Replace('id="3" name="image2.png"', 'id="3" name="image1.png"')
$xml = $xml.Replace('id="0" name="image2.png"', 'id="0" name="image1.png"')
$ftrFirst.Range.WordOpenXML = $xml

# --- Footer (default/primary), holds the Pearson logo currently named "image2.png",
#     docPr/cNvPr id="2" -> rename to "image1.png" ---
$ftrDefault = $sec.Footers(1)
$xml = $ftrDefault.Range.WordOpenXML
$xml = $xml.Replace('id="2" name="image2.png"', 'id="2" name="image1.png"')
$xml = $xml.Replace('id="0" name="image2.png"', 'id="0" name="image1.png"')
$ftrDefault.Range.WordOpenXML = $xml

# --- Header (first-page), holds the BTEC logo currently named "image1.jpg",
#     docPr/cNvPr id="1" -> rename to "image2.jpg" ---
$hdrFirst = $sec.Headers(2)
$xml = $hdrFirst.Range.WordOpenXML
$xml = $xml.Replace('id="1" name="image1.jpg"', 'id="1" name="image2.jpg"')
$xml = $xml.Replace('id="0" name="image1.jpg"', 'id="0" name="image2.jpg"')
$hdrFirst.Range.WordOpenXML = $xml

Write-Output "renamed docPr/cNvPr image names in headers/footers"
